# Generate Report for Handback
# ------------------------------------------------------------------
# This models a "handback" run of the localization-status report:
#   * the Status column flips from "Ready for handoff" to
#     "Handed back: in sync with en-US" on every sheet that shows it
#   * the two data sheets (zh-cn / de-de) get their "Latest Target
#     File" / "Latest Handback File" / "Latest Handback DateTime"
#     columns populated for both rows, with the target-file cell
#     turned into a real hyperlink back to the source .md file
#   * a couple of columns get wider to fit the newly-populated data
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$ovw  = $wb.Worksheets.Item("Overview")
$zh   = $wb.Worksheets.Item("zh-cn")
$de   = $wb.Worksheets.Item("de-de")

# ---- 1. Flip the status text everywhere it appears -----------------
$ovw.Range("E2").Value = $statusNew
$ovw.Range("F2").Value = $statusNew
$ovw.Range("E3").Value = $statusNew
$ovw.Range("F3").Value = $statusNew

$zh.Range("C2").Value = $statusNew
$zh.Range("C3").Value = $statusNew

$de.Range("C2").Value = $statusNew
$de.Range("C3").Value = $statusNew

# ---- 2. Populate handback info for zh-cn ----------------------------
$zhHandbackDateTime = "2016-08-20 00:52:13"

$zh.Hyperlinks.Add(
    $zh.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/594a203fa0a07d476ed2c499042acea1f5d7c28d/e2e/15ef7362-bcf9-4f0e-a1c7-55f0e9fafb6f.md",
    "",
    "",
    "15ef7362-bcf9-4f0e-a1c7-55f0e9fafb6f.md"
) | Out-Null
$zh.Range("J2").Value = "15ef7362-bcf9-4f0e-a1c7-55f0e9fafb6f.80f0f44bc21ac19016f2f33a712f2bcd89283340.zh-cn.xlf"
$zh.Range("K2").Value = $zhHandbackDateTime

$zh.Hyperlinks.Add(
    $zh.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/594a203fa0a07d476ed2c499042acea1f5d7c28d/e2e/8c254ff2-a83e-465d-b8e4-238f5c507734.md",
    "",
    "",
    "8c254ff2-a83e-465d-b8e4-238f5c507734.md"
) | Out-Null
$zh.Range("J3").Value = "8c254ff2-a83e-465d-b8e4-238f5c507734.f6dd222ccea9630494b1a0bb4c7b08fd0460b431.zh-cn.xlf"
$zh.Range("K3").Value = $zhHandbackDateTime

# ---- 3. Populate handback info for de-de ----------------------------
$deHandbackDateTime = "2016-08-20 00:52:19"

$de.Hyperlinks.Add(
    $de.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/594a203fa0a07d476ed2c499042acea1f5d7c28d/e2e/15ef7362-bcf9-4f0e-a1c7-55f0e9fafb6f.md",
    "",
    "",
    "15ef7362-bcf9-4f0e-a1c7-55f0e9fafb6f.md"
) | Out-Null
$de.Range("J2").Value = "15ef7362-bcf9-4f0e-a1c7-55f0e9fafb6f.80f0f44bc21ac19016f2f33a712f2bcd89283340.de-de.xlf"
$de.Range("K2").Value = $deHandbackDateTime

$de.Hyperlinks.Add(
    $de.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/594a203fa0a07d476ed2c499042acea1f5d7c28d/e2e/8c254ff2-a83e-465d-b8e4-238f5c507734.md",
    "",
    "",
    "8c254ff2-a83e-465d-b8e4-238f5c507734.md"
) | Out-Null
$de.Range("J3").Value = "8c254ff2-a83e-465d-b8e4-238f5c507734.f6dd222ccea9630494b1a0bb4c7b08fd0460b431.de-de.xlf"
$de.Range("K3").Value = $deHandbackDateTime

# ---- 4. Widen columns to fit the newly-populated / longer text -----
# Overview: zh-cn / de-de status columns (E, F)
$ovw.Range("E1").ColumnWidth = 29.166666666666668
$ovw.Range("F1").ColumnWidth = 29.166666666666668

# zh-cn / de-de: Status (C) and the new Target/Handback File columns (I, J)
$zh.Range("C1").ColumnWidth = 29.166666666666668
$zh.Range("I1").ColumnWidth = 39.166666666666664
$zh.Range("J1").ColumnWidth = 39.166666666666664

$de.Range("C1").ColumnWidth = 29.166666666666668
$de.Range("I1").ColumnWidth = 39.166666666666664
$de.Range("J1").ColumnWidth = 39.166666666666664

Write-Host "Handback report generated."
